# Apply the IG-publisher metadata refresh:
#  - Jurisdiction value gains the "FRANCE" entry
#  - Date timestamp bumped to the re-run time
# (sheet2 "Include #0" keeps the same visible content; only the
#  underlying shared-string table is reshuffled by Excel on save.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B11").Value = "FRANCE"
$ws.Range("B8").Value = "2025-07-11T12:29:53+00:00"
